# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet, filling them in for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$date = "2012-04-16"
$legislatorName = "曾巨威"
$legislatorId = 1755

# Header row
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Data rows 2 through 19
for ($r = 2; $r -le 19; $r++) {
    # Force text so the ISO-looking date string isn't reinterpreted as a
    # date serial number.
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = $date

    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}

Write-Output "added date/legislator_name/legislator_id columns to 股票 sheet"
